$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 73

$ws.Cells.Item($row, 1).Value = "M2IJF1"
$ws.Cells.Item($row, 2).Value = "Grasa para fusor HP Original"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 150000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 9
$ws.Cells.Item($row, 8).Formula = "=(E73-D73)*G73"
$ws.Cells.Item($row, 9).Formula = "=D73*F73"
$ws.Cells.Item($row, 10).Value = 0
